$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update a handful of odds values on row 7
$ws.Range("G7").Value = 2.3
$ws.Range("I7").Value = 3.2
$ws.Range("Q7").Value = 2.25
$ws.Range("R7").Value = 1.62
$ws.Range("X7").Value = 10
$ws.Range("Y7").Value = 9.5
$ws.Range("Z7").Value = 21
$ws.Range("AH7").Value = 9
$ws.Range("AL7").Value = 29
$ws.Range("AX7").Value = 19
$ws.Range("BB7").Value = 251

# Row 9 (Louisville City vs North Carolina) is removed; the row below it
# (Colorado Springs vs Oakland Roots) shifts up to become the new row 9.
$ws.Rows.Item(9).Delete()
